$wb = $excel.ActiveWorkbook
$docs = $wb.ActiveSheet

# Update selection on the Docs sheet while it is still active.
[void]$docs.Range("A16").Select()

# Insert the new "Plan" sheet (Worksheets.Add inserts before the active sheet,
# so it becomes the first sheet and the new active sheet).
$plan = $wb.Worksheets.Add()
$plan.Name = "Plan"

$plan.Range("A14").Value = "Application:"
$plan.Range("B14").Value = "https://www.rememberthemilk.com/"
$plan.Range("A1").Value = "Setup environment:"
$plan.Range("B2").Value = "IDE (Eclipse/IDEA)"
$plan.Range("B3").Value = "Maven"
$plan.Range("B4").Value = "TestNG"
$plan.Range("B5").Value = "Git"
$plan.Range("B6").Value = "Selenium WebDriver"
$plan.Range("B7").Value = "Windows"
$plan.Range("B8").Value = "Java"
$plan.Range("B10").Value = "Firebug addons for Firefox"
$plan.Range("B11").Value = "FirePath addons for Firefox"
$plan.Range("B9").Value = "Firefox browser"
$plan.Range("B12").Value = "Beyond Compare"

$plan.Hyperlinks.Add($plan.Range("B14"), "https://www.rememberthemilk.com/")
$plan.Range("B14").Style = "Hyperlink"

$plan.Columns.Item(1).ColumnWidth = 14.140625

[void]$plan.Range("B26").Select()
